$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F2: the long "Username : ..." note -----------------------------------
# The embedded "Status Register : 1 : Lanjutkan ke Verifikasi;" line is
# simplified to "Status Register : 1;"
$oldNote = $ws.Range("F2").Text
$newNote = $oldNote -replace "Status Register : 1 : Lanjutkan ke Verifikasi;", "Status Register : 1;"
$ws.Range("F2").Value = $newNote

# --- O2 (STATUS_REGISTER column) -------------------------------------------
# Used to hold the text "1 : Lanjutkan ke Verifikasi"; now it is simply the
# number 1. Align it like the other data cells in the row (left/center,
# wrap text) instead of the old vertical-only alignment.
$o2 = $ws.Range("O2")
$o2.HorizontalAlignment = -4131   # xlLeft
$o2.VerticalAlignment = -4108     # xlCenter
$o2.WrapText = $true
$o2.Value = 1

# --- Row 2 height ------------------------------------------------------
# Shrinks now that O2 no longer needs to wrap the longer text.
$ws.Rows.Item(2).RowHeight = 120
